$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new test case data rows (3 and 4)
$ws.Range("A3").Value = "orch2.dc.local"
$ws.Range("B3").Value = "Default"
$ws.Range("C3").Value = "admin"
$ws.Range("D3").Value = "Password2$"

$ws.Range("A4").Value = "orch2.dc.local"
$ws.Range("B4").Value = "Host"
$ws.Range("C4").Value = "admin"
$ws.Range("D4").Value = "Password2$"

# Update the selected cell to match the saved view state
$ws.Range("B5").Select()
